# Update countries & provincias Spain
# - Sudan overtakes Guam in the ranking (row 168 becomes Sudan, row 169 becomes Guam)
#   with Sudan's stats refreshed; Guam's stats are unchanged, just shifted down a row.
# - Refresh daily stats for several other countries (rows 4, 20, 59, 95, 136, 191).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 695107
$ws.Range("C4").Value = 17537
$ws.Range("D4").Value = 59147
$ws.Range("E4").Value = 599384
$ws.Range("G4").Value = 1959
$ws.Range("H4").Value = 36576

# --- Row 20: Austria ---
$ws.Range("E20").Value = 4451
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 431

# --- Row 59: Moldavia ---
$ws.Range("E59").Value = 1932
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 56

# --- Row 95: Niger ---
$ws.Range("B95").Value = 627
$ws.Range("C95").Value = 43
$ws.Range("D95").Value = 110
$ws.Range("E95").Value = 499
$ws.Range("G95").Value = 4
$ws.Range("H95").Value = 18

# --- Row 136: Somalia ---
$ws.Range("E136").Value = 108
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = 6

# --- Rows 168/169: Sudan moves ahead of Guam ---
# Row 168 was Guam, now becomes Sudan with refreshed totals.
$ws.Range("A168").Value = "Sudan"
$ws.Range("B168").Value = 33
$ws.Range("C168").Value = 1
$ws.Range("D168").Value = 4
$ws.Range("E168").Value = 23
$ws.Range("G168").Value = 1
$ws.Range("H168").Value = 6

# Row 169 was Sudan, now becomes Guam with its (unchanged) totals.
$ws.Range("A169").Value = "Guam"
$ws.Range("B169").Value = 32
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 0
$ws.Range("E169").Value = 31
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 1

# --- Row 191: Granada ---
$ws.Range("D191").Value = 6
$ws.Range("E191").Value = 8
$ws.Range("F191").Value = 4
